# Updates cryptos list data (prices and 1h volume %) per the latest scrape.
# Also fixes three row-order swaps where two coins traded ranking positions
# (rows 19/20, 34/35, 47/48): Coin name + Link are swapped between the pair,
# while Price + Volume(1h) are each set to their own new scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'98.157.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "'3.410.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'256.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "

# Row 6
$ws.Range("D6").Value = "'661.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.86%  "

# Row 7
$ws.Range("E7").Value = "  -4.39%  "

# Row 8
$ws.Range("D8").Value = "'0.436"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.29%  "

# Row 9
$ws.Range("E9").Value = "  -1.30%  "

# Row 10
$ws.Range("E10").Value = "  -0.02%  "

# Row 11
$ws.Range("D11").Value = "'3.407.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.66%  "

# Row 12
$ws.Range("E12").Value = "  +3.44%  "

# Row 13
$ws.Range("D13").Value = "'42.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.86%  "

# Row 14
$ws.Range("D14").Value = "'6.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +15.14%  "

# Row 15
$ws.Range("D15").Value = "'97.804.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "

# Row 16
$ws.Range("E16").Value = "  -1.59%  "

# Row 17
$ws.Range("D17").Value = "'4.044.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.36%  "

# Row 18
$ws.Range("D18").Value = "'9.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +22.17%  "

# Row 19
$ws.Range("B19").Value = "Stellar"
$ws.Range("C19").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D19").Value = "'0.589"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +40.28%  "

# Row 20
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "'3.409.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.26%  "

# Row 21
$ws.Range("D21").Value = "'17.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.68%  "

# Row 22
$ws.Range("D22").Value = "'10.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.12%  "

# Row 23
$ws.Range("D23").Value = "'517.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.78%  "

# Row 24
$ws.Range("D24").Value = "'3.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.17%  "

# Row 25
$ws.Range("E25").Value = "  -3.20%  "

# Row 26
$ws.Range("D26").Value = "'6.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.52%  "

# Row 27
$ws.Range("D27").Value = "'100.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "

# Row 28
$ws.Range("D28").Value = "'12.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "

# Row 29
$ws.Range("D29").Value = "'3.594.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "

# Row 30
$ws.Range("E30").Value = "  +1.12%  "

# Row 31
$ws.Range("E31").Value = "  +7.88%  "

# Row 32
$ws.Range("E32").Value = "  +4.59%  "

# Row 33
$ws.Range("D33").Value = "'0.996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "

# Row 34
$ws.Range("B34").Value = "PolygonEcosystemToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D34").Value = "'0.580"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.87%  "

# Row 35
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("D36").Value = "'2.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.51%  "

# Row 37
$ws.Range("D37").Value = "'30.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.62%  "

# Row 38
$ws.Range("D38").Value = "'7.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.43%  "

# Row 39
$ws.Range("D39").Value = "'1.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.37%  "

# Row 40
$ws.Range("D40").Value = "'534.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.64%  "

# Row 41
$ws.Range("D41").Value = "'0.153"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.35%  "

# Row 42
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("D43").Value = "'0.879"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.50%  "

# Row 44
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("D45").Value = "'9.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +18.86%  "

# Row 46
$ws.Range("D46").Value = "'5.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +20.67%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0430"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.06%  "

# Row 48
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").Value = "'3.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.90%  "

# Row 49
$ws.Range("D49").Value = "'1.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.10%  "

# Row 50
$ws.Range("D50").Value = "'3.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.17%  "

# Row 51
$ws.Range("D51").Value = "'2.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.87%  "
